# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7823
$ws1.Range("F5").Value = 5713
$ws1.Range("F9").Value = 67
$ws1.Range("F11").Value = 274
$ws1.Range("F12").Value = 59

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7823
$ws4.Range("F5").Value = 5713
$ws4.Range("F9").Value = 67
$ws4.Range("F13").Value = 275
$ws4.Range("F14").Value = 59
